$d = $word.ActiveDocument

# The "Vehicular Access to Semi-detached Plot" section embeds a tiny 1x1
# placeholder picture (meant to stand in for a real screenshot). The fix
# swaps that inline picture out for a plain hyperlink run -- styled with
# the built-in "Hyperlink" character style -- whose visible text and
# target both point at the real image hosted on ura.gov.sg.
$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/SD02_Vehicular_Access.jpg?h=100%25&w=100%25"

if ($d.InlineShapes.Count -ge 1) {
    $shp = $d.InlineShapes(1)

    # Keep a handle on the (now-empty) spot the picture used to occupy so
    # we can drop the hyperlink run in exactly the same place.
    $rng = $shp.Range
    $shp.Delete()

    # Hyperlinks.Add(Anchor, Address, SubAddress, ScreenTip, TextToDisplay)
    # inserts a run at $rng showing the URL as text and wires up the
    # external relationship + "Hyperlink" character style automatically.
    $d.Hyperlinks.Add($rng, $url, $null, $null, $url) | Out-Null
}

Write-Host "Replaced picture with hyperlink"
